$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Run Mode" column (C) from "No" to "Yes" for rows 11-28, 30, 31, 33
$rows = @(11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,30,31,33)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = "Yes"
}

# Update sheet view: scroll position and selection
$excel.ActiveWindow.ScrollRow = 32
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C11:C33").Select()
